# Updates NATMI TPM-derived statistics for the F9 -> Lrp1 ligand-receptor pairs sheet.
# The underlying ligand-expressing cell counts were recomputed (TPM re-run), which
# cascades into the derived expression / specificity columns for every row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.249426
$ws.Range("H2").Value = 0.748278
$ws.Range("I2").Value = 0.3330492191161541
$ws.Range("J2").Value = 0.3330492191161541
$ws.Range("M2").Value = 9.172748
$ws.Range("N2").Value = 27.518244
$ws.Range("O2").Value = 0.01445826353606064
$ws.Range("P2").Value = 0.01445826353606064
$ws.Range("Q2").Value = 2.287921842648
$ws.Range("R2").Value = 20.591296583832
$ws.Range("S2").Value = 0.004815313380460561
$ws.Range("T2").Value = 0.00481531338046056

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.249426
$ws.Range("H3").Value = 0.748278
$ws.Range("I3").Value = 0.3330492191161541
$ws.Range("J3").Value = 0.3330492191161541
$ws.Range("O3").Value = 0.2254554169720557
$ws.Range("P3").Value = 0.2254554169720557
$ws.Range("Q3").Value = 35.67678592572
$ws.Range("R3").Value = 321.09107333148
$ws.Range("S3").Value = 0.07508775056805007
$ws.Range("T3").Value = 0.07508775056805006

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.249426
$ws.Range("H4").Value = 0.748278
$ws.Range("I4").Value = 0.3330492191161541
$ws.Range("J4").Value = 0.3330492191161541
$ws.Range("M4").Value = 169.2367096666667
$ws.Range("N4").Value = 507.7101290000001
$ws.Range("O4").Value = 0.2667541884216647
$ws.Range("P4").Value = 0.2667541884216647
$ws.Range("Q4").Value = 42.212035545318
$ws.Range("R4").Value = 379.9083199078621
$ws.Range("S4").Value = 0.08884227414979884
$ws.Range("T4").Value = 0.08884227414979884

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.249426
$ws.Range("H5").Value = 0.748278
$ws.Range("I5").Value = 0.3330492191161541
$ws.Range("J5").Value = 0.3330492191161541
$ws.Range("M5").Value = 16.15031566666667
$ws.Range("N5").Value = 48.450947
$ws.Range("O5").Value = 0.02545644119943506
$ws.Range("P5").Value = 0.02545644119943505
$ws.Range("Q5").Value = 4.028308635474001
$ws.Range("R5").Value = 36.254777719266
$ws.Range("S5").Value = 0.008478247862948138
$ws.Range("T5").Value = 0.008478247862948138

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.249426
$ws.Range("H6").Value = 0.748278
$ws.Range("I6").Value = 0.3330492191161541
$ws.Range("J6").Value = 0.3330492191161541
$ws.Range("M6").Value = 54.744643
$ws.Range("N6").Value = 164.233929
$ws.Range("O6").Value = 0.08628956945961638
$ws.Range("P6").Value = 0.08628956945961638
$ws.Range("Q6").Value = 13.654737324918
$ws.Range("R6").Value = 122.892635924262
$ws.Range("S6").Value = 0.02873867372639437
$ws.Range("T6").Value = 0.02873867372639437

# Row 7
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.249426
$ws.Range("H7").Value = 0.748278
$ws.Range("I7").Value = 0.3330492191161541
$ws.Range("J7").Value = 0.3330492191161541
$ws.Range("M7").Value = 242.0894676666667
$ws.Range("N7").Value = 726.268403
$ws.Range("O7").Value = 0.3815861204111676
$ws.Range("P7").Value = 0.3815861204111676
$ws.Range("Q7").Value = 60.383407562226
$ws.Range("R7").Value = 543.450668060034
$ws.Range("S7").Value = 0.1270869594285021
$ws.Range("T7").Value = 0.1270869594285021

# Row 8
$ws.Range("G8").Value = 0.4994903333333333
$ws.Range("H8").Value = 1.498471
$ws.Range("I8").Value = 0.666950780883846
$ws.Range("J8").Value = 0.666950780883846
$ws.Range("M8").Value = 9.172748
$ws.Range("N8").Value = 27.518244
$ws.Range("O8").Value = 0.01445826353606064
$ws.Range("P8").Value = 0.01445826353606064
$ws.Range("Q8").Value = 4.581698956102667
$ws.Range("R8").Value = 41.235290604924
$ws.Range("S8").Value = 0.009642950155600079
$ws.Range("T8").Value = 0.009642950155600077

# Row 9
$ws.Range("G9").Value = 0.4994903333333333
$ws.Range("H9").Value = 1.498471
$ws.Range("I9").Value = 0.666950780883846
$ws.Range("J9").Value = 0.666950780883846
$ws.Range("O9").Value = 0.2254554169720557
$ws.Range("P9").Value = 0.2254554169720557
$ws.Range("Q9").Value = 71.44487621298445
$ws.Range("R9").Value = 643.00388591686
$ws.Range("S9").Value = 0.1503676664040056
$ws.Range("T9").Value = 0.1503676664040056

# Row 10
$ws.Range("G10").Value = 0.4994903333333333
$ws.Range("H10").Value = 1.498471
$ws.Range("I10").Value = 0.666950780883846
$ws.Range("J10").Value = 0.666950780883846
$ws.Range("M10").Value = 169.2367096666667
$ws.Range("N10").Value = 507.7101290000001
$ws.Range("O10").Value = 0.2667541884216647
$ws.Range("P10").Value = 0.2667541884216647
$ws.Range("Q10").Value = 84.5321005236399
$ws.Range("R10").Value = 760.7889047127591
$ws.Range("S10").Value = 0.1779119142718658
$ws.Range("T10").Value = 0.1779119142718658

# Row 11
$ws.Range("G11").Value = 0.4994903333333333
$ws.Range("H11").Value = 1.498471
$ws.Range("I11").Value = 0.666950780883846
$ws.Range("J11").Value = 0.666950780883846
$ws.Range("M11").Value = 16.15031566666667
$ws.Range("N11").Value = 48.450947
$ws.Range("O11").Value = 0.02545644119943506
$ws.Range("P11").Value = 0.02545644119943505
$ws.Range("Q11").Value = 8.066926555781889
$ws.Range("R11").Value = 72.602339002037
$ws.Range("S11").Value = 0.01697819333648692
$ws.Range("T11").Value = 0.01697819333648692

# Row 12
$ws.Range("G12").Value = 0.4994903333333333
$ws.Range("H12").Value = 1.498471
$ws.Range("I12").Value = 0.666950780883846
$ws.Range("J12").Value = 0.666950780883846
$ws.Range("M12").Value = 54.744643
$ws.Range("N12").Value = 164.233929
$ws.Range("O12").Value = 0.08628956945961638
$ws.Range("P12").Value = 0.08628956945961638
$ws.Range("Q12").Value = 27.34441998028434
$ws.Range("R12").Value = 246.099779822559
$ws.Range("S12").Value = 0.05755089573322202
$ws.Range("T12").Value = 0.05755089573322202

# Row 13
$ws.Range("G13").Value = 0.4994903333333333
$ws.Range("H13").Value = 1.498471
$ws.Range("I13").Value = 0.666950780883846
$ws.Range("J13").Value = 0.666950780883846
$ws.Range("M13").Value = 242.0894676666667
$ws.Range("N13").Value = 726.268403
$ws.Range("O13").Value = 0.3815861204111676
$ws.Range("P13").Value = 0.3815861204111676
$ws.Range("Q13").Value = 120.9213489013126
$ws.Range("R13").Value = 1088.292140111813
$ws.Range("S13").Value = 0.2544991609826655
$ws.Range("T13").Value = 0.2544991609826655
